$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "52.208.93"
Set-TextValue $ws.Range("E2") "  -13.85%  "
Set-TextValue $ws.Range("D3") "2.292.58"
Set-TextValue $ws.Range("E3") "  -21.00%  "
Set-TextValue $ws.Range("E4") "  +0.30%  "
Set-TextValue $ws.Range("D5") "442.39"
Set-TextValue $ws.Range("E5") "  -16.02%  "
Set-TextValue $ws.Range("D6") "119.92"
Set-TextValue $ws.Range("E6") "  -16.18%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.465"
Set-TextValue $ws.Range("E8") "  -15.27%  "
Set-TextValue $ws.Range("D9") "2.294.20"
Set-TextValue $ws.Range("E9") "  -21.03%  "
Set-TextValue $ws.Range("D10") "5.27"
Set-TextValue $ws.Range("E10") "  -11.85%  "
Set-TextValue $ws.Range("D11") "0.0865"
Set-TextValue $ws.Range("E11") "  -19.22%  "
Set-TextValue $ws.Range("D12") "0.301"
Set-TextValue $ws.Range("E12") "  -16.21%  "
Set-TextValue $ws.Range("E13") "  -6.43%  "
Set-TextValue $ws.Range("D14") "52.233.20"
Set-TextValue $ws.Range("E14") "  -13.78%  "
Set-TextValue $ws.Range("D15") "18.70"
Set-TextValue $ws.Range("E15") "  -17.08%  "
Set-TextValue $ws.Range("E16") "  -17.10%  "
Set-TextValue $ws.Range("D17") "2.320.46"
Set-TextValue $ws.Range("E17") "  -20.11%  "
Set-TextValue $ws.Range("D18") "3.94"
Set-TextValue $ws.Range("E18") "  -20.90%  "
Set-TextValue $ws.Range("D19") "296.35"
Set-TextValue $ws.Range("E19") "  -15.50%  "
Set-TextValue $ws.Range("D20") "8.84"
Set-TextValue $ws.Range("E20") "  -23.69%  "
Set-TextValue $ws.Range("D21") "0.999"
Set-TextValue $ws.Range("E21") "  -0.14%  "
Set-TextValue $ws.Range("D22") "5.62"
Set-TextValue $ws.Range("E22") "  -1.68%  "
Set-TextValue $ws.Range("D23") "5.11"
Set-TextValue $ws.Range("E23") "  -21.59%  "
Set-TextValue $ws.Range("D24") "53.35"
Set-TextValue $ws.Range("E24") "  -17.46%  "
Set-TextValue $ws.Range("D25") "0.363"
Set-TextValue $ws.Range("E25") "  -19.68%  "
Set-TextValue $ws.Range("D26") "0.145"
Set-TextValue $ws.Range("E26") "  -18.82%  "
Set-TextValue $ws.Range("D27") "6.87"
Set-TextValue $ws.Range("E27") "  -12.11%  "
Set-TextValue $ws.Range("D28") "0.998"
Set-TextValue $ws.Range("E28") "  -0.10%  "
Set-TextValue $ws.Range("D29") "0.0₃0652"
Set-TextValue $ws.Range("E29") "  -23.21%  "
Set-TextValue $ws.Range("D30") "141.91"
Set-TextValue $ws.Range("E30") "  -5.81%  "
Set-TextValue $ws.Range("D31") "16.71"
Set-TextValue $ws.Range("E31") "  -14.61%  "
Set-TextValue $ws.Range("D32") "1.33"
Set-TextValue $ws.Range("E32") "  -20.47%  "
Set-TextValue $ws.Range("D33") "4.67"
Set-TextValue $ws.Range("E33") "  -16.04%  "
Set-TextValue $ws.Range("D34") "0.813"
Set-TextValue $ws.Range("E34") "  -18.41%  "
Set-TextValue $ws.Range("D35") "3.36"
Set-TextValue $ws.Range("E35") "  -22.00%  "
Set-TextValue $ws.Range("D36") "0.996"
Set-TextValue $ws.Range("E36") "  -0.11%  "
Set-TextValue $ws.Range("D37") "0.985"
Set-TextValue $ws.Range("E37") "  -17.52%  "
Set-TextValue $ws.Range("D38") "31.76"
Set-TextValue $ws.Range("E38") "  -15.68%  "
Set-TextValue $ws.Range("D39") "10.15"
Set-TextValue $ws.Range("E39") "  -1.62%  "
Set-TextValue $ws.Range("D40") "0.555"
Set-TextValue $ws.Range("E40") "  -14.27%  "
Set-TextValue $ws.Range("D41") "0.0504"
Set-TextValue $ws.Range("E41") "  -13.22%  "
Set-TextValue $ws.Range("D42") "3.09"
Set-TextValue $ws.Range("E42") "  -16.91%  "
Set-TextValue $ws.Range("D43") "1.907.41"
Set-TextValue $ws.Range("E43") "  -16.76%  "
Set-TextValue $ws.Range("D44") "1.16"
Set-TextValue $ws.Range("E44") "  -20.98%  "
Set-TextValue $ws.Range("D45") "0.0817"
Set-TextValue $ws.Range("E45") "  -11.14%  "
Set-TextValue $ws.Range("E46") "  -14.31%  "
Set-TextValue $ws.Range("D47") "4.14"
Set-TextValue $ws.Range("E47") "  -16.19%  "
Set-TextValue $ws.Range("D48") "15.48"
Set-TextValue $ws.Range("E48") "  -24.32%  "
Set-TextValue $ws.Range("E49") "  -5.13%  "
Set-TextValue $ws.Range("D50") "4.44"
Set-TextValue $ws.Range("E50") "  -12.97%  "
Set-TextValue $ws.Range("D51") "14.85"
Set-TextValue $ws.Range("E51") "  -18.61%  "
